$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Name = "mc constants"
